$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($addr in @("D5", "D6", "D9", "D11", "D12", "D13", "D14", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D33", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.203.09"
$ws.Range("E2").Value = "  +3.06%  "
$ws.Range("D3").Value = "3.054.86"
$ws.Range("E3").Value = "  +2.97%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "598.24"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").Value = "152.88"
$ws.Range("E6").Value = "  +7.85%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.050.66"
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  +11.58%  "
$ws.Range("D11").Value = "0.154"
$ws.Range("E11").Value = "  +7.85%  "
$ws.Range("D12").Value = "0.464"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").Value = "0.0000237"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("D14").Value = "35.27"
$ws.Range("E14").Value = "  +4.26%  "
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "3.557.79"
$ws.Range("E16").Value = "  +2.85%  "
$ws.Range("D17").Value = "63.088.03"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "7.08"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "3.042.60"
$ws.Range("E19").Value = "  +2.65%  "
$ws.Range("D20").Value = "456.23"
$ws.Range("E20").Value = "  +2.23%  "
$ws.Range("D21").Value = "14.34"
$ws.Range("E21").Value = "  +3.60%  "
$ws.Range("D22").Value = "0.698"
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +3.11%  "
$ws.Range("D24").Value = "82.92"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  +6.67%  "
$ws.Range("D26").Value = "10.97"
$ws.Range("E26").Value = "  +11.17%  "
$ws.Range("D27").Value = "12.29"
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("E29").Value = "  +3.76%  "
$ws.Range("D30").Value = "7.46"
$ws.Range("E30").Value = "  +9.39%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +6.10%  "
$ws.Range("D33").Value = "27.80"
$ws.Range("E33").Value = "  +3.26%  "
$ws.Range("E34").Value = "  +5.60%  "
$ws.Range("D35").Value = "0.0₃0872"
$ws.Range("E35").Value = "  +12.14%  "
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("D37").Value = "5.93"
$ws.Range("E37").Value = "  +3.70%  "
$ws.Range("D38").Value = "3.17"
$ws.Range("E38").Value = "  +15.04%  "
$ws.Range("D39").Value = "2.13"
$ws.Range("E39").Value = "  +3.69%  "
$ws.Range("D40").Value = "50.64"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("D41").Value = "9.13"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "0.125"
$ws.Range("E42").Value = "  +4.55%  "
$ws.Range("D43").Value = "0.297"
$ws.Range("E43").Value = "  +13.23%  "
$ws.Range("D44").Value = "41.41"
$ws.Range("E44").Value = "  +11.82%  "
$ws.Range("D45").Value = "396.94"
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("D46").Value = "0.0358"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").Value = "2.757.00"
$ws.Range("E47").Value = "  +2.75%  "
$ws.Range("D48").Value = "132.23"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D50").Value = "2.22"
$ws.Range("E50").Value = "  +3.77%  "
$ws.Range("D51").Value = "24.25"
$ws.Range("E51").Value = "  +4.33%  "
